# IRYO-vaccination_data.xlsx update: add the 2021-05-21 (Fri) daily row and
# refresh the running totals / "as of" note, per the source workbook's daily
# refresh pattern (new row inserted right under the cumulative-total row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 5 (pushes the old row 5 -> row 6, etc.)
$ws.Rows.Item(5).Insert()

# 2) Copy number formats/styles down from the row that used to be row 5
#    (now row 6) into the freshly inserted row 5, so it matches the rest
#    of the daily-data rows instead of picking up blank defaults.
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill the new row 5 with the 2021-05-21 (Fri) daily figures.
$ws.Range("A5").Value = 44337
$ws.Range("B5").Value = "(金)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 99918
$ws.Range("E5").Value = 149103

# 4) Update the cumulative-total row (row 4) to the new running totals.
$ws.Range("D4").Value = 3965411
$ws.Range("E4").Value = 2472976

# 5) Refresh the "as of" note from 5/20 to 5/21.
$ws.Range("E2").Value = "（5月21日時点）"
